$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 8606911.488
$ws.Range("L2").Value = 9695297.536
$ws.Range("M2").Value = 9425941.504000001
$ws.Range("N2").Value = 10745706.496
$ws.Range("K3").Value = 8533017.088
$ws.Range("L3").Value = 9619281.92
$ws.Range("M3").Value = 9344338.944
$ws.Range("N3").Value = 10657928.192
$ws.Range("K4").Value = 143548
$ws.Range("L4").Value = 246919.008
$ws.Range("M4").Value = 207664.992
$ws.Range("N4").Value = 287188
$ws.Range("K5").Value = 8366720
$ws.Range("L5").Value = 9337278.464
$ws.Range("M5").Value = 9096536.063999999
$ws.Range("N5").Value = 10343438.336
$ws.Range("K9").Value = 0
$ws.Range("M9").Value = 40138
$ws.Range("N9").Value = 0
$ws.Range("K11").Value = 22749
$ws.Range("L11").Value = 35085
$ws.Range("N11").Value = 27302
$ws.Range("K12").Value = 12280
$ws.Range("L12").Value = 13142
$ws.Range("M12").Value = 18875
$ws.Range("N12").Value = 25229
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("K19").Value = 12280
$ws.Range("L19").Value = 13142
$ws.Range("M19").Value = 18875
$ws.Range("N19").Value = 25229
$ws.Range("K23").Value = 45865
$ws.Range("L23").Value = 46259
$ws.Range("M23").Value = 45158
$ws.Range("N23").Value = 44063
$ws.Range("K24").Value = 15749
$ws.Range("L24").Value = 16615
$ws.Range("M24").Value = 17569
$ws.Range("N24").Value = 18487
$ws.Range("K26").Value = 8606911.488
$ws.Range("L26").Value = 9695297.536
$ws.Range("M26").Value = 9425941.504000001
$ws.Range("N26").Value = 10745706.496
$ws.Range("K27").Value = 7814678.016
$ws.Range("L27").Value = 8888539.136
$ws.Range("M27").Value = 8649956.352
$ws.Range("N27").Value = 9924889.6
$ws.Range("K29").Value = 2623
$ws.Range("L29").Value = 2712
$ws.Range("M29").Value = 2740
$ws.Range("N29").Value = 4727
$ws.Range("K30").Value = 91720
$ws.Range("L30").Value = 111652
$ws.Range("M30").Value = 123217
$ws.Range("N30").Value = 120825
$ws.Range("K31").Value = 7480474.112
$ws.Range("L31").Value = 8254886.912
$ws.Range("M31").Value = 8150203.904
$ws.Range("N31").Value = 9319444.48
$ws.Range("K34").Value = 229642
$ws.Range("L34").Value = 500476.992
$ws.Range("M34").Value = 344720
$ws.Range("N34").Value = 431670.016
$ws.Range("K35").Value = 10219
$ws.Range("L35").Value = 18811
$ws.Range("M35").Value = 29075
$ws.Range("N35").Value = 48224
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 0
$ws.Range("N37").Value = 0
$ws.Range("K47").Value = 792233.024
$ws.Range("L47").Value = 806758.976
$ws.Range("M47").Value = 775985.024
$ws.Range("N47").Value = 820817.024
$ws.Range("K48").Value = 674940.032
$ws.Range("L48").Value = 674940.032
$ws.Range("M48").Value = 674940.032
$ws.Range("N48").Value = 674940.032
$ws.Range("K49").Value = -30193
$ws.Range("L49").Value = -30193
$ws.Range("M49").Value = -30193
$ws.Range("N49").Value = -30193
$ws.Range("K51").Value = 119729
$ws.Range("L51").Value = 119729
$ws.Range("M51").Value = 50431
$ws.Range("N51").Value = 180316.992
$ws.Range("K52").Value = 33123
$ws.Range("L52").Value = 46628
$ws.Range("M52").Value = 86778
$ws.Range("K53").Value = -4824
$ws.Range("L53").Value = -3866
$ws.Range("M53").Value = -5222
$ws.Range("N53").Value = -3610
$ws.Range("K54").Value = -542
$ws.Range("L54").Value = -479
$ws.Range("M54").Value = -749
$ws.Range("N54").Value = -637
$ws.Range("K57").ClearContents()
$ws.Range("L57").ClearContents()
$ws.Range("M57").ClearContents()
$ws.Range("N57").ClearContents()
$ws.Range("K58").ClearContents()
$ws.Range("L58").ClearContents()
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()
$ws.Range("K59").Value = 47424
$ws.Range("L59").Value = 56978
$ws.Range("M59").Value = 64862
$ws.Range("N59").Value = 77400.992
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = 0
$ws.Range("N60").Value = 0
$ws.Range("K61").Value = 47424
$ws.Range("L61").Value = 56978
$ws.Range("M61").Value = 64862
$ws.Range("N61").Value = 77400.992
$ws.Range("K63").Value = -49860
$ws.Range("L63").Value = -50221
$ws.Range("M63").Value = -52913
$ws.Range("N63").Value = -83053.008
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 0
$ws.Range("N64").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
$ws.Range("N65").Value = 0
$ws.Range("K66").Value = 13
$ws.Range("L66").Value = 21
$ws.Range("M66").Value = -43
$ws.Range("N66").Value = 472
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = 0
$ws.Range("N67").Value = 0
$ws.Range("K68").Value = 54083
$ws.Range("L68").Value = 46172
$ws.Range("M68").Value = 41983
$ws.Range("N68").Value = 46910
$ws.Range("K69").Value = 1020760
$ws.Range("L69").Value = 1276967.04
$ws.Range("M69").Value = 1418843.008
$ws.Range("N69").Value = 2593467.904
$ws.Range("K70").Value = -966676.992
$ws.Range("L70").Value = -1230795.008
$ws.Range("M70").Value = -1376860.032
$ws.Range("N70").Value = -2546557.952
$ws.Range("K71").ClearContents()
$ws.Range("L71").ClearContents()
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("K74").Value = 51660
$ws.Range("L74").Value = 52950
$ws.Range("M74").Value = 53889
$ws.Range("N74").Value = 41729.992
$ws.Range("K75").Value = -4198
$ws.Range("L75").Value = -15534
$ws.Range("M75").Value = -37310
$ws.Range("N75").Value = 4513
$ws.Range("K76").Value = -14339
$ws.Range("L76").Value = 1288
$ws.Range("M76").Value = 23571
$ws.Range("N76").Value = -3136
$ws.Range("K77").ClearContents()
$ws.Range("L77").ClearContents()
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = 0
$ws.Range("N79").Value = 0
$ws.Range("K80").Value = 33123
$ws.Range("L80").Value = 38704
$ws.Range("M80").Value = 40150
$ws.Range("N80").Value = 43107
